$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 397.75555  # H17: 1057.7556 -> 397.75555
$ws.Cells.Item(17, 10).Value = 397.75555  # J17: 1057.7556 -> 397.75555
$ws.Cells.Item(17, 12).Value = 1193.26665  # L17: 3173.2668 -> 1193.26665
$ws.Cells.Item(17, 14).Value = -1529.26665  # N17: -3509.2668 -> -1529.26665
$ws.Cells.Item(40, 8).Value = 1701.9429  # H40: 1707.7222 -> 1701.9429
$ws.Cells.Item(40, 9).Value = 1591.8846  # I40: 1571.96 -> 1591.8846
$ws.Cells.Item(40, 10).Value = 2019.8889  # J40: 2016.2727 -> 2019.8889
$ws.Cells.Item(40, 11).Value = 1591.8846  # K40: 1571.96 -> 1591.8846
$ws.Cells.Item(40, 12).Value = 2019.8889  # L40: 2016.2727 -> 2019.8889
$ws.Cells.Item(40, 13).Value = -1416.8846  # M40: -1396.96 -> -1416.8846
$ws.Cells.Item(40, 14).Value = -2369.8889  # N40: -2366.2727 -> -2369.8889
$ws.Cells.Item(86, 8).Value = 1800  # H86: 1778.0741 -> 1800
$ws.Cells.Item(86, 9).Value = 1800.0714  # I86: 1775 -> 1800.0714
$ws.Cells.Item(86, 10).Value = 1799.5  # J86: 1782.5454 -> 1799.5
$ws.Cells.Item(86, 11).Value = 1800.0714  # K86: 1775 -> 1800.0714
$ws.Cells.Item(86, 12).Value = 1799.5  # L86: 1782.5454 -> 1799.5
$ws.Cells.Item(86, 13).Value = -677.0714  # M86: -652 -> -677.0714
$ws.Cells.Item(86, 14).Value = -4045.5  # N86: -4028.5454 -> -4045.5
$ws.Cells.Item(89, 8).Value = 1800  # H89: 1778.0741 -> 1800
$ws.Cells.Item(89, 9).Value = 1800.0714  # I89: 1775 -> 1800.0714
$ws.Cells.Item(89, 10).Value = 1799.5  # J89: 1782.5454 -> 1799.5
$ws.Cells.Item(89, 11).Value = 9000.357  # K89: 8875 -> 9000.357
$ws.Cells.Item(89, 12).Value = 8997.5  # L89: 8912.726999999999 -> 8997.5
$ws.Cells.Item(89, 13).Value = -3384.357  # M89: -3259 -> -3384.357
$ws.Cells.Item(89, 14).Value = -20229.5  # N89: -20144.727 -> -20229.5
$ws.Cells.Item(92, 8).Value = 1838  # H92: 2442.5 -> 1838
$ws.Cells.Item(92, 9).Value = 409.8  # I92: 575 -> 409.8
$ws.Cells.Item(92, 10).Value = 4218.3335  # J92: 6177.5 -> 4218.3335
$ws.Cells.Item(92, 11).Value = 409.8  # K92: 575 -> 409.8
$ws.Cells.Item(92, 12).Value = 4218.3335  # L92: 6177.5 -> 4218.3335
$ws.Cells.Item(92, 13).Value = 838.2  # M92: 673 -> 838.2
$ws.Cells.Item(92, 14).Value = -6714.3335  # N92: -8673.5 -> -6714.3335
$ws.Cells.Item(100, 8).Value = 1455.4  # H100: 1379.6666 -> 1455.4
$ws.Cells.Item(100, 9).Value = 1392.3334  # I100: 1548.5 -> 1392.3334
$ws.Cells.Item(100, 10).Value = 1550  # J100: 1331.4286 -> 1550
$ws.Cells.Item(100, 11).Value = 1392.3334  # K100: 1548.5 -> 1392.3334
$ws.Cells.Item(100, 12).Value = 1550  # L100: 1331.4286 -> 1550
$ws.Cells.Item(100, 13).Value = -851.3334  # M100: -1007.5 -> -851.3334
$ws.Cells.Item(100, 14).Value = -2632  # N100: -2413.4286 -> -2632
$ws.Cells.Item(101, 8).Value = 1172.8334  # H101: 1171.7273 -> 1172.8334
$ws.Cells.Item(101, 10).Value = 2194  # J101: 2446.25 -> 2194
$ws.Cells.Item(101, 12).Value = 6582  # L101: 7338.75 -> 6582
$ws.Cells.Item(101, 14).Value = -9826  # N101: -10582.75 -> -9826
$ws.Cells.Item(125, 8).Value = 2651.818  # H125: 2514.1667 -> 2651.818
$ws.Cells.Item(125, 9).Value = 1178  # I125: 1133.5 -> 1178
$ws.Cells.Item(125, 11).Value = 10602  # K125: 10201.5 -> 10602
$ws.Cells.Item(125, 13).Value = -8142  # M125: -7741.5 -> -8142
$ws.Cells.Item(137, 8).Value = 34017.613  # H137: 36274.1 -> 34017.613
$ws.Cells.Item(137, 10).Value = 93149.55  # J137: 113560.89 -> 93149.55
$ws.Cells.Item(137, 12).Value = 279448.65  # L137: 340682.67 -> 279448.65
$ws.Cells.Item(137, 14).Value = -284548.65  # N137: -345782.67 -> -284548.65

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(92, 8).Value = 21525  # H92: 24699.666 -> 21525
$ws.Cells.Item(92, 10).Value = 21525  # J92: 24699.666 -> 21525
$ws.Cells.Item(92, 12).Value = 21525  # L92: 24699.666 -> 21525
$ws.Cells.Item(92, 14).Value = -26517  # N92: -29691.666 -> -26517
$ws.Cells.Item(97, 8).Value = 809.875  # H97: 557.3158 -> 809.875
$ws.Cells.Item(97, 9).Value = 714.1429000000001  # I97: 522.4375 -> 714.1429000000001
$ws.Cells.Item(97, 10).Value = 1480  # J97: 743.3333 -> 1480
$ws.Cells.Item(97, 11).Value = 714.1429000000001  # K97: 522.4375 -> 714.1429000000001
$ws.Cells.Item(97, 12).Value = 1480  # L97: 743.3333 -> 1480
$ws.Cells.Item(97, 13).Value = -218.1429000000001  # M97: -26.4375 -> -218.1429000000001
$ws.Cells.Item(97, 14).Value = -2472  # N97: -1735.3333 -> -2472
$ws.Cells.Item(110, 8).Value = 2909.2917  # H110: 2812.92 -> 2909.2917
$ws.Cells.Item(110, 9).Value = 2528.5881  # I110: 2415.889 -> 2528.5881
$ws.Cells.Item(110, 11).Value = 2528.5881  # K110: 2415.889 -> 2528.5881
$ws.Cells.Item(110, 13).Value = -483.5880999999999  # M110: -370.8890000000001 -> -483.5880999999999
$ws.Cells.Item(123, 8).Value = 42377.8  # H123: 59089.5 -> 42377.8
$ws.Cells.Item(123, 10).Value = 42377.8  # J123: 59089.5 -> 42377.8
$ws.Cells.Item(123, 12).Value = 42377.8  # L123: 59089.5 -> 42377.8
$ws.Cells.Item(123, 14).Value = -52177.8  # N123: -68889.5 -> -52177.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 466.625  # H64: 429.29413 -> 466.625
$ws.Cells.Item(64, 9).Value = 652.5  # I64: 507.125 -> 652.5
$ws.Cells.Item(64, 10).Value = 355.1  # J64: 360.1111 -> 355.1
$ws.Cells.Item(64, 11).Value = 652.5  # K64: 507.125 -> 652.5
$ws.Cells.Item(64, 12).Value = 355.1  # L64: 360.1111 -> 355.1
$ws.Cells.Item(64, 13).Value = -427.5  # M64: -282.125 -> -427.5
$ws.Cells.Item(64, 14).Value = -805.1  # N64: -810.1111000000001 -> -805.1
$ws.Cells.Item(67, 8).Value = 466.625  # H67: 429.29413 -> 466.625
$ws.Cells.Item(67, 9).Value = 652.5  # I67: 507.125 -> 652.5
$ws.Cells.Item(67, 10).Value = 355.1  # J67: 360.1111 -> 355.1
$ws.Cells.Item(67, 11).Value = 652.5  # K67: 507.125 -> 652.5
$ws.Cells.Item(67, 12).Value = 355.1  # L67: 360.1111 -> 355.1
$ws.Cells.Item(67, 13).Value = 127.5  # M67: 272.875 -> 127.5
$ws.Cells.Item(67, 14).Value = -1915.1  # N67: -1920.1111 -> -1915.1
$ws.Cells.Item(94, 8).Value = 1437.375  # H94: 1544.1428 -> 1437.375
$ws.Cells.Item(94, 9).Value = 1285.5714  # I94: 1384.8334 -> 1285.5714
$ws.Cells.Item(94, 11).Value = 1285.5714  # K94: 1384.8334 -> 1285.5714
$ws.Cells.Item(94, 13).Value = -834.5714  # M94: -933.8334 -> -834.5714
$ws.Cells.Item(99, 8).Value = 3192.7856  # H99: 3120.6667 -> 3192.7856
$ws.Cells.Item(99, 9).Value = 2847.5  # I99: 2680 -> 2847.5
$ws.Cells.Item(99, 10).Value = 3330.9  # J99: 3341 -> 3330.9
$ws.Cells.Item(99, 11).Value = 2847.5  # K99: 2680 -> 2847.5
$ws.Cells.Item(99, 12).Value = 3330.9  # L99: 3341 -> 3330.9
$ws.Cells.Item(99, 13).Value = -1349.5  # M99: -1182 -> -1349.5
$ws.Cells.Item(99, 14).Value = -6326.9  # N99: -6337 -> -6326.9

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2560.1667  # H132: 1686.6 -> 2560.1667
$ws.Cells.Item(132, 9).Value = 1332.1111  # I132: 1012.44446 -> 1332.1111
$ws.Cells.Item(132, 10).Value = 3788.2222  # J132: 3961.875 -> 3788.2222
$ws.Cells.Item(132, 11).Value = 3996.3333  # K132: 3037.33338 -> 3996.3333
$ws.Cells.Item(132, 12).Value = 11364.6666  # L132: 11885.625 -> 11364.6666
$ws.Cells.Item(132, 13).Value = -1466.3333  # M132: -507.33338 -> -1466.3333
$ws.Cells.Item(132, 14).Value = -16424.6666  # N132: -16945.625 -> -16424.6666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 53227.633  # H4: 51536.418 -> 53227.633
$ws.Cells.Item(4, 10).Value = 66425.336  # J4: 63800.32 -> 66425.336
$ws.Cells.Item(4, 12).Value = 199276.008  # L4: 191400.96 -> 199276.008
$ws.Cells.Item(4, 14).Value = -199500.008  # N4: -191624.96 -> -199500.008
$ws.Cells.Item(56, 8).Value = 9125  # H56: 6666.6665 -> 9125
$ws.Cells.Item(56, 9).Value = 9125  # I56: 6666.6665 -> 9125
$ws.Cells.Item(56, 11).Value = 9125  # K56: 6666.6665 -> 9125
$ws.Cells.Item(56, 13).Value = -8595  # M56: -6136.6665 -> -8595
$ws.Cells.Item(131, 8).Value = 739.16364  # H131: 766.25 -> 739.16364
$ws.Cells.Item(131, 9).Value = 338  # I131: 340.5263 -> 338
$ws.Cells.Item(131, 10).Value = 968.4  # J131: 946 -> 968.4
$ws.Cells.Item(131, 11).Value = 1014  # K131: 1021.5789 -> 1014
$ws.Cells.Item(131, 12).Value = 2905.2  # L131: 2838 -> 2905.2
$ws.Cells.Item(131, 13).Value = 4026  # M131: 4018.4211 -> 4026
$ws.Cells.Item(131, 14).Value = -12985.2  # N131: -12918 -> -12985.2
$ws.Cells.Item(136, 8).Value = 4418.75  # H136: 4477.273 -> 4418.75
$ws.Cells.Item(136, 9).Value = 2116.6667  # I136: 2062.5 -> 2116.6667
$ws.Cells.Item(136, 10).Value = 5800  # J136: 5857.143 -> 5800
$ws.Cells.Item(136, 11).Value = 6350.000100000001  # K136: 6187.5 -> 6350.000100000001
$ws.Cells.Item(136, 12).Value = 17400  # L136: 17571.429 -> 17400
$ws.Cells.Item(136, 13).Value = -1250.000100000001  # M136: -1087.5 -> -1250.000100000001
$ws.Cells.Item(136, 14).Value = -27600  # N136: -27771.429 -> -27600
$ws.Cells.Item(138, 8).Value = 2930  # H138: 2724.56 -> 2930
$ws.Cells.Item(138, 9).Value = 0  # I138: 1632.5 -> 0
$ws.Cells.Item(138, 10).Value = 2930  # J138: 2932.5715 -> 2930
$ws.Cells.Item(138, 11).Value = 0  # K138: 4897.5 -> 0
$ws.Cells.Item(138, 12).Value = 8790  # L138: 8797.7145 -> 8790
$ws.Cells.Item(138, 13).ClearContents()  # M138: was 242.5
$ws.Cells.Item(138, 14).Value = -19070  # N138: -19077.7145 -> -19070

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 435635.6  # H97: 1157.85 -> 435635.6
$ws.Cells.Item(97, 9).Value = 556494.25  # I97: 1321.9286 -> 556494.25
$ws.Cells.Item(97, 10).Value = 544.4  # J97: 775 -> 544.4
$ws.Cells.Item(97, 11).Value = 556494.25  # K97: 1321.9286 -> 556494.25
$ws.Cells.Item(97, 12).Value = 544.4  # L97: 775 -> 544.4
$ws.Cells.Item(97, 13).Value = -555998.25  # M97: -825.9286 -> -555998.25
$ws.Cells.Item(97, 14).Value = -1536.4  # N97: -1767 -> -1536.4
$ws.Cells.Item(132, 8).Value = 4807.2  # H132: 5116.913 -> 4807.2
$ws.Cells.Item(132, 9).Value = 5570.3076  # I132: 5910.3335 -> 5570.3076
$ws.Cells.Item(132, 10).Value = 3980.5  # J132: 4251.364 -> 3980.5
$ws.Cells.Item(132, 11).Value = 16710.9228  # K132: 17731.0005 -> 16710.9228
$ws.Cells.Item(132, 12).Value = 11941.5  # L132: 12754.092 -> 11941.5
$ws.Cells.Item(132, 13).Value = -14180.9228  # M132: -15201.0005 -> -14180.9228
$ws.Cells.Item(132, 14).Value = -17001.5  # N132: -17814.092 -> -17001.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 435298.34  # H22: 417164.25 -> 435298.34
$ws.Cells.Item(22, 9).Value = 588568.8  # I22: 555875 -> 588568.8
$ws.Cells.Item(22, 11).Value = 588568.8  # K22: 555875 -> 588568.8
$ws.Cells.Item(22, 13).Value = -588273.8  # M22: -555580 -> -588273.8
$ws.Cells.Item(27, 8).Value = 435298.34  # H27: 417164.25 -> 435298.34
$ws.Cells.Item(27, 9).Value = 588568.8  # I27: 555875 -> 588568.8
$ws.Cells.Item(27, 11).Value = 588568.8  # K27: 555875 -> 588568.8
$ws.Cells.Item(27, 13).Value = -588461.8  # M27: -555768 -> -588461.8
$ws.Cells.Item(93, 8).Value = 16107  # H93: 18766.334 -> 16107
$ws.Cells.Item(93, 9).Value = 18674.834  # I93: 18766.334 -> 18674.834
$ws.Cells.Item(93, 10).Value = 700  # J93: 0 -> 700
$ws.Cells.Item(93, 11).Value = 18674.834  # K93: 18766.334 -> 18674.834
$ws.Cells.Item(93, 12).Value = 700  # L93: 0 -> 700
$ws.Cells.Item(93, 13).Value = -17426.834  # M93: -17518.334 -> -17426.834
$ws.Cells.Item(93, 14).Value = -3196  # N93: None -> -3196
$ws.Cells.Item(100, 8).Value = 35716388  # H100: 166670600 -> 35716388
$ws.Cells.Item(100, 9).Value = 2600.375  # I100: 5900 -> 2600.375
$ws.Cells.Item(100, 10).Value = 83334776  # J100: 500000000 -> 83334776
$ws.Cells.Item(100, 11).Value = 2600.375  # K100: 5900 -> 2600.375
$ws.Cells.Item(100, 12).Value = 83334776  # L100: 500000000 -> 83334776
$ws.Cells.Item(100, 13).Value = -2059.375  # M100: -5359 -> -2059.375
$ws.Cells.Item(100, 14).Value = -83335858  # N100: -500001082 -> -83335858
$ws.Cells.Item(108, 8).Value = 0  # H108: 30000 -> 0
$ws.Cells.Item(108, 10).Value = 0  # J108: 30000 -> 0
$ws.Cells.Item(108, 12).Value = 0  # L108: 30000 -> 0
$ws.Cells.Item(108, 14).ClearContents()  # N108: was -37680
$ws.Cells.Item(122, 8).Value = 2599.5  # H122: 0 -> 2599.5
$ws.Cells.Item(122, 10).Value = 2599.5  # J122: 0 -> 2599.5
$ws.Cells.Item(122, 12).Value = 7798.5  # L122: 0 -> 7798.5
$ws.Cells.Item(122, 14).Value = -12698.5  # N122: None -> -12698.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 568.8461  # H113: 683.6667 -> 568.8461
$ws.Cells.Item(113, 9).Value = 467  # I113: 477.2857 -> 467
$ws.Cells.Item(113, 10).Value = 731.8  # J113: 972.6 -> 731.8
$ws.Cells.Item(113, 11).Value = 1401  # K113: 1431.8571 -> 1401
$ws.Cells.Item(113, 12).Value = 2195.4  # L113: 2917.8 -> 2195.4
$ws.Cells.Item(113, 13).Value = 769  # M113: 738.1428999999998 -> 769
$ws.Cells.Item(113, 14).Value = -6535.4  # N113: -7257.8 -> -6535.4
$ws.Cells.Item(136, 8).Value = 2263.0833  # H136: 2519.762 -> 2263.0833
$ws.Cells.Item(136, 9).Value = 1350.7142  # I136: 1591.909 -> 1350.7142
$ws.Cells.Item(136, 11).Value = 4052.1426  # K136: 4775.727000000001 -> 4052.1426
$ws.Cells.Item(136, 13).Value = -1502.1426  # M136: -2225.727000000001 -> -1502.1426
$ws.Cells.Item(138, 8).Value = 38970.25  # H138: 39922.668 -> 38970.25
$ws.Cells.Item(138, 10).Value = 38970.25  # J138: 39922.668 -> 38970.25
$ws.Cells.Item(138, 12).Value = 38970.25  # L138: 39922.668 -> 38970.25
$ws.Cells.Item(138, 14).Value = -49250.25  # N138: -50202.668 -> -49250.25
